# Update countries & provincias Spain
# Refresh of the COVID-19 "Pais" data snapshot: the source table got a new
# pull (06:52 -> 07:22) which re-sorted a handful of countries whose totals
# changed rank, plus updated the case counters for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp banner (row 1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 07:22"

# --- Re-sorted block around San Marino / Mauricio / Nigeria / Kirguistan --
$ws.Range("A104").Value = "Kirguistan"
$ws.Range("B104").Value = 339
$ws.Range("C104").Value = 41
$ws.Range("D104").Value = 44
$ws.Range("E104").Value = 290
$ws.Range("F104").Value = 5
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 5

$ws.Range("A105").Value = "Mauricio"
$ws.Range("B105").Value = 318
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 23
$ws.Range("E105").Value = 286
$ws.Range("F105").Value = 3
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 9

$ws.Range("A106").Value = "Nigeria"
$ws.Range("B106").Value = 305
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 58
$ws.Range("E106").Value = 240
$ws.Range("F106").Value = 2
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 7

# --- Re-sorted block around Bahamas / Guyana / Zambia / Puerto Rico / Liberia
$ws.Range("A151").Value = "Guyana"
$ws.Range("B151").Value = 40
$ws.Range("C151").Value = 3
$ws.Range("D151").Value = 8
$ws.Range("E151").Value = 26
$ws.Range("F151").Value = 3
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 6

$ws.Range("A152").Value = "Zambia"
$ws.Range("B152").Value = 40
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 25
$ws.Range("E152").Value = 13
$ws.Range("F152").Value = 1
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 2

$ws.Range("A153").Value = "Puerto Rico"
$ws.Range("B153").Value = 39
$ws.Range("C153").Value = 0
$ws.Range("D153").Value = 1
$ws.Range("E153").Value = 36
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 2

$ws.Range("A154").Value = "Liberia"
$ws.Range("B154").Value = 37
$ws.Range("C154").Value = 0
$ws.Range("D154").Value = 3
$ws.Range("E154").Value = 29
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 5

# --- Swap of Somalia / Antigua y Barbuda -----------------------------------
$ws.Range("A164").Value = "Antigua y Barbuda"
$ws.Range("B164").Value = 21
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 0
$ws.Range("E164").Value = 19
$ws.Range("F164").Value = 1
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 2

$ws.Range("A165").Value = "Somalia"
$ws.Range("B165").Value = 21
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 1
$ws.Range("E165").Value = 19
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 1
